$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 22
$startA = 10002
$startB = 110021
$rowCount = 9

for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $startA + $i
    $ws.Cells.Item($r, 2).Value = $startB + $i
    $ws.Cells.Item($r, 3).Value = "eng"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Select the rows below the data, mirroring the author's post-paste selection.
[void]$ws.Range("A31:XFD1048576").Select()

# Add page setup info present in the final worksheet (portrait orientation).
$ws.PageSetup.Orientation = 1
